# Applies the "promise" document updates described by the commit:
#   upate title module promise, roundgarbage, roundmoney
#
# Word Find.Execute positional signature used below:
#   FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#   MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace
#
# wdFindContinue = 1 ; wdReplaceOne = 1 ; wdReplaceAll = 2

function Replace-Text {
    param([string]$Find, [string]$Replace)
    $d.Content.Find.Execute($Find, $true, $false, $false, $false, $false, $true, 1, $false, $Replace, 2)
}

$d = $word.ActiveDocument

# 1. Document number (เลขที่) : 121212121212 -> 12000
Replace-Text "121212121212" "12000"

# 2 & 3. Contract start date "23 มิถุนายน 2562" (appears twice) -> "2 ก.ค. 2562"
Replace-Text "23 มิถุนายน 2562" "2 ก.ค. 2562"

# 4. Contract end date "23 มิถุนายน 2563" -> "30 ก.ย. 2562"
Replace-Text "23 มิถุนายน 2563" "30 ก.ย. 2562"

# 5. Payment period "รายเดือน" -> "รายครั้ง"
Replace-Text "รายเดือน" "รายครั้ง"

# 6 & 7. Monthly rate "1000" and its Thai text "หนึ่งพันบาทถ้วน" cleared out
Replace-Text "เดือนละ 1000บาท (หนึ่งพันบาทถ้วน" "เดือนละ บาท ("

# 8. Pickup count "จัดเก็บ3ครั้งต่อ" -> "จัดเก็บ0ครั้งต่อ"
Replace-Text "จัดเก็บ3ครั้งต่อ" "จัดเก็บ0ครั้งต่อ"

# 9 & 10. Annual total "10000" -> "20000", Thai text "หนึ่งหมื่นบาทถ้วน" cleared out
Replace-Text "ต่อปี 10000 บาท (หนึ่งหมื่นบาทถ้วน" "ต่อปี 20000 บาท ("
